$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows (and create new rows) with new A/B values per the commit diff
$ws.Range("A2").Value = "RODZICE  Wykłady dla rodziców IMOPEKSIS Goczałkowice/Nieporęt  Spojrzenie dziecka na udział w treningach sportowych"
$ws.Range("B2").Value = "https://www.youtube.com/watch?v=i_wLW6hwob4"
$ws.Range("A3").Value = "RODZICE  Wykłady dla rodziców IMOPEKSIS Goczałkowice/Nieporęt  Uczestnictwo w sporcie oczami dziecka a oczami dorosłego"
$ws.Range("B3").Value = "https://www.youtube.com/watch?v=vMh8ukTNUYo"
$ws.Range("A4").Value = "RODZICE  Szkolenie dla rodziców – Imopeksis Suchy Las 08.12.2023  Szkolenie dla rodziców – cz.1"
$ws.Range("B4").Value = "https://www.youtube.com/watch?v=l5rPMwXT0qc"
$ws.Range("A5").Value = "RODZICE  Szkolenie dla rodziców – Imopeksis Suchy Las 08.12.2023  Szkolenie dla rodziców – cz.2"
$ws.Range("B5").Value = "https://www.youtube.com/watch?v=R5BxAirxUYk"
$ws.Range("A6").Value = "RODZICE  Szkolenie dla rodziców – Imopeksis cz.3  Rodzice cz.1"
$ws.Range("B6").Value = "https://www.youtube.com/watch?v=-UIq3Luu6go"
$ws.Range("A7").Value = "RODZICE  Szkolenie dla rodziców – Imopeksis cz.3  Rodzice cz.2"
$ws.Range("B7").Value = "https://www.youtube.com/watch?v=yrmY5qzzse8"
$ws.Range("A8").Value = "RODZICE  Szkolenie dla rodziców – Imopeksis cz.3  Rodzice cz.3"
$ws.Range("B8").Value = "https://www.youtube.com/watch?v=279q6pYHmvA"
$ws.Range("A9").Value = "RODZICE  Szkolenie dla rodziców – Imopeksis cz.3  Rodzice cz.4"
$ws.Range("B9").Value = "https://www.youtube.com/watch?v=vK8mOmlf5jc"
$ws.Range("A10").Value = "RODZICE  Szkolenia dla rodziców – Imopeksis cz. 4  Rozumienie kompetencji psychospołecznych u dzieci i młodzieży cz.1"
$ws.Range("B10").Value = "https://www.youtube.com/watch?v=l_RWPZHMDjs"
$ws.Range("A11").Value = "RODZICE  Szkolenia dla rodziców – Imopeksis cz. 4  Rozumienie kompetencji psychospołecznych u dzieci i młodzieży cz.2"
$ws.Range("B11").Value = "https://www.youtube.com/watch?v=994Y8-ltBFI"
$ws.Range("A12").Value = "RODZICE  Szkolenia dla rodziców – Imopeksis cz. 4  Imopeksis Gdańsk 22.02.2024 cz1"
$ws.Range("B12").Value = "https://www.youtube.com/watch?v=SBBr9Gc_zBk"
$ws.Range("A13").Value = "RODZICE  Szkolenia dla rodziców – Imopeksis cz. 4  Klinika Imopeksis Wykłady cz 1"
$ws.Range("B13").Value = "https://www.youtube.com/watch?v=sKBReSoDqhQ"
$ws.Range("A14").Value = "RODZICE  Szkolenia dla rodziców – Imopeksis cz. 5  Imopeksis OZPN Nyska 7.02.2024 HCF dla rodziców"
$ws.Range("B14").Value = "https://www.youtube.com/watch?v=-1VVuFuYNac"
$ws.Range("A15").Value = "RODZICE  Szkolenia dla rodziców – Imopeksis cz. 5  Imopeksis Łomża 18.02.2024 Szkolenie dla rodziców cz2"
$ws.Range("B15").Value = "https://www.youtube.com/watch?v=cQZUoecBfjw"
$ws.Range("A16").Value = "RODZICE  Szkolenia dla rodziców – Imopeksis cz. 5  Platforma temat 3 cz. 1"
$ws.Range("B16").Value = "https://www.youtube.com/watch?v=gWUjBZwlkos"
$ws.Range("A17").Value = "RODZICE  Szkolenia dla rodziców – Imopeksis cz. 5  Platforma temat 3 cz. 2"
$ws.Range("B17").Value = "https://www.youtube.com/watch?v=yQSYLbxf9Ao"
$ws.Range("A18").Value = "TRENERZY  Szkolenie dla trenerów – Klinika IMOPEKSIS  Poruszanie się bez piłki"
$ws.Range("B18").Value = "https://www.youtube.com/watch?v=dJ5NcBptShQ"
$ws.Range("A19").Value = "TRENERZY  Szkolenie dla trenerów – Klinika IMOPEKSIS  Rzut z miejsca"
$ws.Range("B19").Value = "https://www.youtube.com/watch?v=84HpdfJnRfI"
$ws.Range("A20").Value = "TRENERZY  Szkolenie dla trenerów – Klinika IMOPEKSIS  Podania specjalne"
$ws.Range("B20").Value = "https://www.youtube.com/watch?v=yxJ5mg6uiq4"
$ws.Range("A21").Value = "TRENERZY  Szkolenie dla trenerów – Klinika IMOPEKSIS cz.2  Podsumowanie Treningu 02.07.2023"
$ws.Range("B21").Value = "https://www.youtube.com/watch?v=DezQxekRLoQ"
$ws.Range("A22").Value = "TRENERZY  Szkolenie dla trenerów – Klinika IMOPEKSIS cz.2  Periodyzacja Techniczna"
$ws.Range("B22").Value = "https://www.youtube.com/watch?v=NUqFVpwPyjs"
$ws.Range("A23").Value = "TRENERZY  Szkolenie dla trenerów cz.4  Szkolenie 1"
$ws.Range("B23").Value = "https://www.youtube.com/watch?v=5zkgjd5aQhM"
$ws.Range("A24").Value = "TRENERZY  Szkolenie dla trenerów cz.4  Szkolenie 2"
$ws.Range("B24").Value = "https://www.youtube.com/watch?v=QqopwhMZEvg"
$ws.Range("A25").Value = "TRENERZY  Szkolenie dla trenerów cz.4  Szkolenie 3"
$ws.Range("B25").Value = "https://www.youtube.com/watch?v=tynSnrho_to"
$ws.Range("A26").Value = "TRENERZY  Szkolenie dla trenerów cz.4  Szkolenie 4"
$ws.Range("A26").Font.Name = "Arial"
$ws.Range("A26").Font.Size = 14
$ws.Range("B26").Value = "https://www.youtube.com/watch?v=-fiih0HvyJI"
$ws.Range("A27").Value = "TRENERZY  Szkolenie dla trenerów cz.5  Imopeksis – profesor Tadeusz Huciński cz5"
$ws.Range("A27").Font.Name = "Arial"
$ws.Range("A27").Font.Size = 14
$ws.Range("B27").Value = "https://www.youtube.com/watch?v=gcFTL83vCAY"
$ws.Range("A28").Value = "TRENERZY  Szkolenie dla trenerów cz.5  Klinika Imopeksis Wykłady trenerów cz 1"
$ws.Range("B28").Value = "https://www.youtube.com/watch?v=sKBReSoDqhQ"
$ws.Range("A29").Value = "TRENERZY  Szkolenie dla trenerów cz.6  Klinika Imopeksis Wykłady trenerów cz 2"
$ws.Range("B29").Value = "https://www.youtube.com/watch?v=qjwGnzLZ2Xg"
$ws.Range("A32").Value = "Nauczyciele Szkolenie dla nauczycieli cz.3  Szkolenie 1"
$ws.Range("B32").Value = "https://www.youtube.com/watch?v=paHp-y3y6l8"
$ws.Range("A33").Value = "Nauczyciele Szkolenie dla nauczycieli cz.3  Szkolenie 2"
$ws.Range("B33").Value = "https://www.youtube.com/watch?v=geHr8JuekDs"
$ws.Range("A34").Value = "Nauczyciele Szkolenie dla nauczycieli cz.3  Szkolenie 3"
$ws.Range("B34").Value = "https://www.youtube.com/watch?v=vIUz2Ia_r0c"
$ws.Range("A35").Value = "Nauczyciele Szkolenie dla nauczycieli cz.3  Szkolenie 4"
$ws.Range("B35").Value = "https://www.youtube.com/watch?v=969fZWWqMgw"
$ws.Range("A36").Value = "Nauczyciele Szkolenie dla nauczycieli cz.4  Imopeksis Radom 19.02.2024 wychowanie fizyczne cz.1"
$ws.Range("B36").Value = "https://www.youtube.com/watch?v=2xbQuoOWtJM"
$ws.Range("A37").Value = "Nauczyciele Szkolenie dla nauczycieli cz.4  Imopeksis Radom 19.02.2024 wychowanie fizyczne cz.2"
$ws.Range("B37").Value = "https://www.youtube.com/watch?v=JAZ4x6xWFCI"
$ws.Range("A38").Value = "Nauczyciele Szkolenie dla nauczycieli cz.4  Imopeksis Łódź 14.01.2024 cz1"
$ws.Range("B38").Value = "https://www.youtube.com/watch?v=YN_TnPeFPGs"
$ws.Range("A39").Value = "Nauczyciele Szkolenie dla nauczycieli cz.4  Imopeksis SP Spalona 7.02.2024 cz 2"
$ws.Range("B39").Value = "https://www.youtube.com/watch?v=-LtRlAf-w2U"
$ws.Range("A40").Value = "Nauczyciele Szkolenie dla nauczycieli cz.5  Imopeksis Kurs trenera koszykówki Łódź cz3 9.03.2024"
$ws.Range("B40").Value = "https://www.youtube.com/watch?v=9ldhbZgqX3w"
$ws.Range("A41").Value = "Nauczyciele Szkolenie dla nauczycieli cz.5  Imopeksis 28.02.2024 Wodzisław Śląski cz2"
$ws.Range("B41").Value = "https://www.youtube.com/watch?v=ApYOp4CW2TU"
$ws.Range("A42").Value = "Nauczyciele Szkolenie dla nauczycieli cz.5  Platforma temat 3 cz. 1"
$ws.Range("B42").Value = "https://www.youtube.com/watch?v=gWUjBZwlkos"
$ws.Range("A43").Value = "Nauczyciele Szkolenie dla nauczycieli cz.5  Imopeksis – profesor Tadeusz Huciński cz6"
$ws.Range("B43").Value = "https://www.youtube.com/watch?v=QgxrlE3hlbg"
$ws.Range("A44").Value = "Nauczyciele Szkolenie dla nauczycieli cz.6  Imopeksis 21.03.2024 Kwidzyn rada pedagogiczna"
$ws.Range("A44").Font.Name = "Arial"
$ws.Range("A44").Font.Size = 14
$ws.Range("B44").Value = "https://www.youtube.com/watch?v=qISNAeFY2Xs"
$ws.Range("A45").Value = "Nauczyciele Szkolenie dla nauczycieli cz.6  Imopeksis 18.04.2024 Włocławek cz2"
$ws.Range("A45").Font.Name = "Arial"
$ws.Range("A45").Font.Size = 14
$ws.Range("B45").Value = "https://www.youtube.com/watch?v=SeCIhc1pWF8"
$ws.Range("A46").Value = "Nauczyciele Szkolenie dla nauczycieli cz.6  Imopeksis Gdańsk 22.02.2024 cz2"
$ws.Range("B46").Value = "https://www.youtube.com/watch?v=zTfv14-LbQY"
$ws.Range("A47").Value = "Nauczyciele Szkolenie dla nauczycieli cz.6  Imopeksis Toruń 28.12.2023 Kurs cz2"
$ws.Range("B47").Value = "https://www.youtube.com/watch?v=5fsyLsMxKCQ"
$ws.Range("A48").Value = "Nauczyciele Szkolenie dla nauczycieli cz.7  Imopeksis Poświętne 15.01.2024 cz2"
$ws.Range("B48").Value = "https://www.youtube.com/watch?v=XOpH0_MyTBw"
$ws.Range("A49").Value = "Nauczyciele Szkolenie dla nauczycieli cz.7  Imopeksis Radom 19.02.2024 wychowanie fizyczne cz1"
$ws.Range("B49").Value = "https://www.youtube.com/watch?v=2xbQuoOWtJM"
$ws.Range("A50").Value = "Nauczyciele  Szkolenie dla nauczycieli cz.7  Imopeksis Łomża 18.02.2024 Szkolenie dla rodziców cz1"
$ws.Range("B50").Value = "https://www.youtube.com/watch?v=HDiTigcOHUE"
$ws.Range("A51").Value = "Nauczyciele Szkolenie dla nauczycieli cz.7  Imopeksis Gdańsk 22.02.2024 cz3"
$ws.Range("B51").Value = "https://www.youtube.com/watch?v=jPFuFqWaR0w"
$ws.Range("A52").Value = "Nauczyciele Szkolenie dla nauczycieli cz.8  Imopeksis 11.06.2024 Ciechanów cz2"
$ws.Range("B52").Value = "https://www.youtube.com/watch?v=nU6ySJNysIg"
$ws.Range("A53").Value = "Nauczyciele Szkolenie dla nauczycieli cz.8  12.05.2024 Imopeksis Kielce cz1"
$ws.Range("B53").Value = "https://www.youtube.com/watch?v=zZGZn5Epmw8"
$ws.Range("A54").Value = " Nauczyciele  Szkolenie dla nauczycieli cz.8  Imopeksis 10.06.2024 Wińsko cz1"
$ws.Range("B54").Value = "https://www.youtube.com/watch?v=t3ulhcF5ATs"
$ws.Range("A55").Value = "Nauczyciele Szkolenie dla nauczycieli cz.8  Imopeksis Tychy 1-3.03.2024 Wykład 1"
$ws.Range("B55").Value = "https://www.youtube.com/watch?v=5fa5JrNRAKI"

# Clear rows that no longer have data (removed from the list)
$ws.Range("A30:B31").Clear()
$ws.Range("A56:B57").Clear()

# Restore selection to match the committed state
$ws.Range("A18:XFD19").Select()
